# Auto-generated: applies the scheduled-runner value updates to the
# Ultima_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1666.6666
$ws.Range("J32").Value = 1666.6666
$ws.Range("L32").Value = 1666.6666
$ws.Range("N32").Value = -2318.6666

$ws.Range("H70").Value = 1461.6
$ws.Range("J70").Value = 1645.3
$ws.Range("L70").Value = 4935.9
$ws.Range("N70").Value = -5475.9

$ws.Range("H73").Value = 1461.6
$ws.Range("J73").Value = 1645.3
$ws.Range("L73").Value = 4935.9
$ws.Range("N73").Value = -6807.9

$ws.Range("H138").Value = 2254.8057
$ws.Range("I138").Value = 1740.1666
$ws.Range("J138").Value = 2769.4443
$ws.Range("K138").Value = 5220.4998
$ws.Range("L138").Value = 8308.332900000001
$ws.Range("M138").Value = -80.4997999999996
$ws.Range("N138").Value = -18588.3329

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 14950
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 14950
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 14950
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -15180

$ws.Range("H74").Value = 20837372
$ws.Range("I74").Value = 33335566
$ws.Range("K74").Value = 33335566
$ws.Range("M74").Value = -33334692

$ws.Range("H77").Value = 20837372
$ws.Range("I77").Value = 33335566
$ws.Range("K77").Value = 166677830
$ws.Range("M77").Value = -166673462

$ws.Range("H133").Value = 58904.2
$ws.Range("J133").Value = 58904.2
$ws.Range("L133").Value = 58904.2
$ws.Range("N133").Value = -63964.2

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1398.1428
$ws.Range("I64").Value = 1972
$ws.Range("J64").Value = 633
$ws.Range("K64").Value = 1972
$ws.Range("L64").Value = 633
$ws.Range("M64").Value = -1747
$ws.Range("N64").Value = -1083

$ws.Range("H67").Value = 1398.1428
$ws.Range("I67").Value = 1972
$ws.Range("J67").Value = 633
$ws.Range("K67").Value = 1972
$ws.Range("L67").Value = 633
$ws.Range("M67").Value = -1192
$ws.Range("N67").Value = -2193

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 72584.71
$ws.Range("I16").Value = 112052.445
$ws.Range("J16").Value = 1542.8
$ws.Range("K16").Value = 112052.445
$ws.Range("L16").Value = 1542.8
$ws.Range("M16").Value = -111765.445
$ws.Range("N16").Value = -2116.8

$ws.Range("H31").Value = 9528410
$ws.Range("I31").Value = 6416.696
$ws.Range("K31").Value = 6416.696
$ws.Range("M31").Value = -6121.696

$ws.Range("H34").Value = 9528410
$ws.Range("I34").Value = 6416.696
$ws.Range("K34").Value = 6416.696
$ws.Range("M34").Value = -6214.696

$ws.Range("H58").Value = 2779.1333
$ws.Range("I58").Value = 947.8333
$ws.Range("K58").Value = 947.8333
$ws.Range("M58").Value = -744.8333

$ws.Range("H99").Value = 1828.2727
$ws.Range("I99").Value = 1159.4
$ws.Range("J99").Value = 2385.6667
$ws.Range("K99").Value = 1159.4
$ws.Range("L99").Value = 2385.6667
$ws.Range("M99").Value = 338.5999999999999
$ws.Range("N99").Value = -5381.6667

$ws.Range("H113").Value = 72584.71
$ws.Range("I113").Value = 112052.445
$ws.Range("J113").Value = 1542.8
$ws.Range("K113").Value = 112052.445
$ws.Range("L113").Value = 1542.8
$ws.Range("M113").Value = -109882.445
$ws.Range("N113").Value = -5882.8

$ws.Range("H122").Value = 1726.6
$ws.Range("I122").Value = 1783.2
$ws.Range("J122").Value = 1556.8
$ws.Range("K122").Value = 5349.6
$ws.Range("L122").Value = 4670.4
$ws.Range("M122").Value = -2899.6
$ws.Range("N122").Value = -9570.4

$ws.Range("H126").Value = 1828.2727
$ws.Range("I126").Value = 1159.4
$ws.Range("J126").Value = 2385.6667
$ws.Range("K126").Value = 3478.2
$ws.Range("L126").Value = 7157.000100000001
$ws.Range("M126").Value = -1008.2
$ws.Range("N126").Value = -12097.0001

$ws.Range("H136").Value = 2779.1333
$ws.Range("I136").Value = 947.8333
$ws.Range("K136").Value = 2843.4999
$ws.Range("M136").Value = -293.4998999999998

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3833.25
$ws.Range("J88").Value = 3833.25
$ws.Range("L88").Value = 11499.75
$ws.Range("N88").Value = -12355.75

$ws.Range("H91").Value = 3833.25
$ws.Range("J91").Value = 3833.25
$ws.Range("L91").Value = 11499.75
$ws.Range("N91").Value = -14463.75

$ws.Range("H97").Value = 7998.4
$ws.Range("I97").Value = 34164.332
$ws.Range("J97").Value = 1456.9166
$ws.Range("K97").Value = 102492.996
$ws.Range("L97").Value = 4370.7498
$ws.Range("M97").Value = -101996.996
$ws.Range("N97").Value = -5362.7498

$ws.Range("H132").Value = 973.75
$ws.Range("I132").Value = 493.33334
$ws.Range("J132").Value = 1133.8889
$ws.Range("K132").Value = 4440.00006
$ws.Range("L132").Value = 10205.0001
$ws.Range("M132").Value = -1910.00006
$ws.Range("N132").Value = -15265.0001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4992
$ws.Range("I5").Value = 4992
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4992
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4880
$ws.Range("N5").Value = $null

$ws.Range("H46").Value = 8000
$ws.Range("I46").Value = 8000
$ws.Range("K46").Value = 8000
$ws.Range("M46").Value = -7844

$ws.Range("H122").Value = 1794.5555
$ws.Range("I122").Value = 1732.75
$ws.Range("J122").Value = 1918.1666
$ws.Range("K122").Value = 5198.25
$ws.Range("L122").Value = 5754.4998
$ws.Range("M122").Value = -2748.25
$ws.Range("N122").Value = -10654.4998

$ws.Range("H138").Value = 58349.332
$ws.Range("J138").Value = 58349.332
$ws.Range("L138").Value = 58349.332
$ws.Range("N138").Value = -68629.332

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10000750

$ws.Range("H22").Value = 878.9375
$ws.Range("I22").Value = 466.25
$ws.Range("J22").Value = 1291.625
$ws.Range("K22").Value = 466.25
$ws.Range("L22").Value = 1291.625
$ws.Range("M22").Value = -171.25
$ws.Range("N22").Value = -1881.625

$ws.Range("H27").Value = 878.9375
$ws.Range("I27").Value = 466.25
$ws.Range("J27").Value = 1291.625
$ws.Range("K27").Value = 466.25
$ws.Range("L27").Value = 1291.625
$ws.Range("M27").Value = -359.25
$ws.Range("N27").Value = -1505.625

$ws.Range("H40").Value = 6189.1055
$ws.Range("I40").Value = 7787.5557
$ws.Range("J40").Value = 4750.5
$ws.Range("K40").Value = 7787.5557
$ws.Range("L40").Value = 4750.5
$ws.Range("M40").Value = -7651.5557
$ws.Range("N40").Value = -5022.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2188.0334
$ws.Range("I122").Value = 2264
$ws.Range("J122").Value = 2074.0833
$ws.Range("K122").Value = 6792
$ws.Range("L122").Value = 6222.249899999999
$ws.Range("M122").Value = -4342
$ws.Range("N122").Value = -11122.2499

$ws.Range("H126").Value = 3419.6897
$ws.Range("I126").Value = 1716.6364
$ws.Range("J126").Value = 8772.143
$ws.Range("K126").Value = 5149.9092
$ws.Range("L126").Value = 26316.429
$ws.Range("M126").Value = -2679.9092
$ws.Range("N126").Value = -31256.429

